$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel date serial for 18 Jan 2018 (matches the "date_test" column elsewhere on the sheet)
$testDateSerial = 43118

# --- Row 9 (values only, formats already correct: A/B = style 5, C = style 9) ---
$ws.Range("A9").Value = "wiki_search_v1_4"
$ws.Range("B9").Value = "wikipedia_search"
$ws.Range("C9").Value = $testDateSerial
$ws.Range("D9").Value = "la fonction donne une page Wikipédia correspondant à l'entité"
$ws.Range("E9").Value = "ok"

# --- Row 10 (only D/E are new; A/B/C stay blank, unchanged) ---
$ws.Range("D10").Value = "la fonction ne renvoie rien si le mot n'existe pas"
$ws.Range("E10").Value = "ok"

# --- Row 11 (only D/E are new) ---
$ws.Range("D11").Value = "la fonction ne renvoie rien si le mot est mal orthographié"
$ws.Range("E11").Value = "ok"

# --- Rows 12-19: C switches from style 4 (blank) to style 9 (date), and the
#     former style-4 formatting moves onto the new D column entry. Grab the
#     formats from existing cells that already carry those styles (C3 = date
#     style 9, C20 = blank style 4) and paste-special them onto the targets
#     so we reuse the existing style records instead of synthesizing new ones. ---

# Row 12
$ws.Range("A12").Value = "g6_polarity_feel_v1_2"
$ws.Range("B12").Value = "load_dict"
$ws.Range("C12").Value = $testDateSerial
$ws.Range("D12").Value = "On ne fait pas de tests sur cette fonction"

# Row 13
$ws.Range("B13").Value = "df_to_dict"
$ws.Range("C13").Value = $testDateSerial
$ws.Range("D13").Value = "On ne fait pas de tests sur cette fonction"

# Row 14
$ws.Range("B14").Value = "Tokenize"
$ws.Range("C14").Value = $testDateSerial
$ws.Range("D14").Value = "On doit retourner une liste ayant  dans chacune des cases un seul mot /caractère (et sans espaces)  correspondant au texte sur lequel on applique cette fonction"
$ws.Range("E14").Value = "OK"

# Row 15
$ws.Range("B15").Value = "del_stop_word_list"
$ws.Range("C15").Value = $testDateSerial
$ws.Range("D15").Value = "La fonction doit retourner une liste de mots  ne contenant plus aucun stopwords"
$ws.Range("E15").Value = "OK"

# Row 16
$ws.Range("B16").Value = "List_to_text"
$ws.Range("C16").Value = $testDateSerial
$ws.Range("D16").Value = "La fonction doit transformer une liste de mots en un seul string"
$ws.Range("E16").Value = "OK"

# Row 17
$ws.Range("B17").Value = "Text_to_list"
$ws.Range("C17").Value = $testDateSerial
$ws.Range("D17").Value = "La fonction doit transformer un string en une liste de mots sans garder les stopwords."
$ws.Range("E17").Value = "OK"

# Row 18
$ws.Range("B18").Value = "feel_polarity"
$ws.Range("C18").Value = $testDateSerial
$ws.Range("D18").Value = "On ne fait pas de tests sur cette fonction"

# Row 19
$ws.Range("B19").Value = "feel_polarity_main"
$ws.Range("C19").Value = $testDateSerial
$ws.Range("D19").Value = "On vérifie que pour une liste donnée, le dictionnaire renvoyé contient bien la moyenne de chacune des composantes de notre vecteur somme ( vecteur qui contient polarité et les émotions générales d'une liste donnée)"
$ws.Range("E19").Value = "OK"

# Re-apply the correct number format (date, style 9) to C12:C19 by copying
# the format from C3, which already carries that style.
$ws.Range("C3").Copy()
$ws.Range("C12:C19").PasteSpecial(-4122)

# Re-apply the formatting the D column used to have on C12:C19 (style 4,
# vertical-center / blank numeric format) onto the newly-populated D cells,
# using C20 (still style 4, untouched by this edit) as the donor.
$ws.Range("C20").Copy()
$ws.Range("D12:D19").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Update the view state saved with the sheet (scroll position reset, new
# active selection).
$ws.Range("C22").Select()
